$d = $word.ActiveDocument

# Paragraph 1: update date in header line
$d.Content.Find.Execute('⚡️🚀המאמר היומי של מייק 07.08.24: ⚡️🚀', $true, $false, $false, $false, $false, $true, 1, $false, '⚡️🚀המאמר היומי של מייק 06.08.24: ⚡️🚀', 2) | Out-Null

# Paragraph 2: replace title text and add a line break
$d.Content.Find.Execute('Language Model Can Listen While Speaking', $true, $false, $false, $false, $false, $true, 1, $false, 'TurboEdit: Text-Based Image Editing Using Few-Step Diffusion Models', 2) | Out-Null
$p2 = $d.Paragraphs(2)
$p2.Range.InsertAfter([char]11)

# Paragraphs 3 and 4: replace body text
$d.Content.Find.Execute(' המאמר שמשך את תשומת ליבי בגלל שמו הקליט. המאמר מציע ארכיטקטורה של מודל Speech Language Model או SLM שיודע להקשיב תוך כדי שהוא מדבר, כלומר מודל full duplex (מושג מתחום התקשורת). בדרך כלל ל- SLM יש שני משטר עבודה: הקשבה או דיבור, כלומר המודל או מדבר או מקשיב. המאמר מעשיר את מרחב היכולות של SLM ומצייד אותו ביכולת להקשיב תוך כדי שהוא מדבר. מעניין שהמודל גם יכול לעצור אם הוא מזהה שיש דיבור (לא רעש) ומגיב עליו (בדיבור) לאחר מכן. ', $true, $false, $false, $false, $false, $true, 1, $false, 'חוזרים לסקור מאמרים על מודלי דיפוזיה עם מאמר כחול לבן של קבוצת חוקרים מאוניברסיטת תל אביב. הם מציעים שיטה מעניינת לעריכה מהירה של תמונה. כלומר בהינתן תמונה עם פרומפט נתון c אנו רוצים ליצור תמונה עם פרומפט אחר c1.', 2) | Out-Null
$d.Content.Find.Execute('הארכיטקטורה של המודל המוצע LSLM מורכב מרכיבים סטנדרטיים. יש מודל שקולט אות דיבור, מחלק אותו לטוקנים (האות במקטעי זמן שונים) מקודד אותו לוקטור אמבדינג ומאזין אותו לדקודר. תפקיד הדקודר הוא לקחת בחשבון את ייצוג של טוקני הדיבור שנקלטו קודם וגם ייצוג טוקני הדיבור שנוצרו על ידי המודל כדי ליצור את הפלט הבא (אות הדיבור) של המודל. ֿכאמור לפעמים הדקודר מחליט שהוא צריך לעבור למצב האזנה ולפעמים הוא צריך לעבור למצב הדיבור.', $true, $false, $false, $false, $false, $true, 1, $false, 'כמו שאתם זוכרים מודלי דיפוזיה מגנרטים תמונה על ידי הסרה רעש הדרגתית (denoising). בכל שלב המודל חוזה כמה רעש צריך להסיר מהתמונה והרעש המשוערך הזה מחוסר מהתמונה המורעשת באיטרציה הקודמת. השיטה הפשוטה לעשות עריכה של תמונה היא:', 2) | Out-Null

# Paragraph 5: replace body text
$d.Content.Find.Execute('כלומר הדקודר במקרה הזה הוא vocoder המקבל כקלט את אות הדיבור הנקלט בנוסף לאות הדיבור המגונרט על ה-vocoder עצמו לפני. ', $true, $false, $false, $false, $false, $true, 1, $false, 'להחסיר מהתמונה(המקורית) באיטרציה t את הרעש הזה המשוערך עם פרומפט c (כמו שעושים כאשר אין עריכה) ', 2) | Out-Null

# Paragraph 6: replace URL paragraph with new text, then append new paragraphs
$d.Content.Find.Execute('https://arxiv.org/pdf/2408.02622', $true, $false, $false, $false, $false, $true, 1, $false, 'להוסיף אל התוצאה את התוחלת המשוערכת של התמונה המורעשת(הערוכה) עם הפרומפט c1 החדש (עם התמונה המורעשת הערוכה. ', 2) | Out-Null

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.Text = 'כלומר בכל איטרציה מתקנים את הסרת הרעש בכיוון הפרומפט החדש.'
$lastPara.Range.InsertParagraphAfter()
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.Text = 'דרך אגב ניתן שערוך הרעש הנוסף באיטרציה t ושערוך תוחלת התמונה אחרי הסרת הרעש אלו שתי בעיות שקולות, כלומר אחת מהן היא פשוט רפרמטריזציה של השנייה מבחינת השערוך.'
$lastPara.Range.InsertParagraphAfter()
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.Text = 'הבעיה בשיטה הפשוטה לעריכת תמונות שהיא לא עובדת טוב ויוצרת ארטיפקטים בתמונה הערוכה. המחברים מנצלים מחקר קודם שמצא שהסקייל של הרעש (כלומר ההפרש בין התמונה המורעשת לתוחלתה) לא מתנהג לפי הסקייל של התהליך הקדמי של הדיפוזיה של התמונה המקורית (שבו מוסיפים רעש עם שונות עולה לתמונה עד שזו הופכת לרעש טהור). הרעש שנוצר במהלך עריכה כזו הוא בעל שונות משמעות גדולה יותר מאשר זה של התמונה המקורית.'
$lastPara.Range.InsertParagraphAfter()
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.Text = 'אז המחברים מציעים להחסיר מהתמונה המורעשת המקורית באיטרציה t את שערוך התוחלת של התמונה המורעשת עבור האיטרציה t+d עבור d חיובי שהם מצאו. כלומר לוקחים תמונה x_t ומזינים אותה למודל שערוך התוחלת עם מספר איטרציה t+d. בסוף מכוונים את התמונה עם שערוך תוחלת המשוערכת של התמונה הערוכה עם איטרציה t+d.'
$lastPara.Range.InsertParagraphAfter()
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.Text = 'בנוסף המאמר מציע דרך מעניינת לווסת את ״עוצמת העריכה״ בצורה דומה ל classifier guidance כדי לכוון את התוצאה של מודל דיפוזיה גנרטיבי ללא פרומפט עבור פרומפט נתון. הפעם על ידי ניתוח של נוסחת העריכה המחברים משקול של מרחק cross-prompt (הפרש שערוך התוחלת עבור התמונה הערוכה המורעשת עבור פרומפטים c ו- c1) לבין מרחק cross-trajectory שמודד הפרש בין חיזוי התוחלת בין התמונה הרגילה לתמונה המשוערכת). משקול כזה מאפשר לבצע את העריכה בפחות איטרציות denoising.'
$lastPara.Range.InsertParagraphAfter()
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.Text = 'מאמר כתוב יפה ובהחלט מומלץ'
$lastPara.Range.InsertParagraphAfter()
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.Text = 'https://arxiv.org/abs/2408.00735'

